# Apply the "merges main and updates historic data" edit:
#  - On the "attribute" sheet, remove the historic "mort" (row 10) and
#    "actualCountID" (row 14) attribute rows, shifting everything below
#    up.
#  - Leave "attribute" as the active sheet/tab (instead of
#    "code_definitions"), with row 13 selected as the active row after
#    the edits.

$wb = $excel.ActiveWorkbook

$wsAttr = $wb.Worksheets.Item("attribute")
$wsCode = $wb.Worksheets.Item("code_definitions")

# Delete from the bottom up so row numbers of earlier rows don't shift
# out from under us.
$wsAttr.Rows.Item(14).Delete()
$wsAttr.Rows.Item(10).Delete()

# Make "attribute" the active sheet/tab.
$wsAttr.Activate()

# Select row 13 (entire row) as the active cell/selection on "attribute".
$wsAttr.Rows.Item(13).Select()

# Keep "code_definitions" selection where it was (B18), just no longer
# the active tab.
$wsCode.Range("B18").Select()

$wsAttr.Activate()
